$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 166.85
$ws.Range("I39").Value = 80.888885
$ws.Range("J39").Value = 237.18182
$ws.Range("K39").Value = 242.666655
$ws.Range("L39").Value = 711.5454599999999
$ws.Range("M39").Value = 53.33334500000001
$ws.Range("N39").Value = -1303.54546
$ws.Range("H40").Value = 2788.652
$ws.Range("I40").Value = 1717.7
$ws.Range("J40").Value = 3612.4614
$ws.Range("K40").Value = 1717.7
$ws.Range("L40").Value = 3612.4614
$ws.Range("M40").Value = -1542.7
$ws.Range("N40").Value = -3962.4614
$ws.Range("H57").Value = 99780
$ws.Range("J57").Value = 99780
$ws.Range("L57").Value = 299340
$ws.Range("N57").Value = -300338
$ws.Range("H64").Value = 7665.4546
$ws.Range("J64").Value = 8340.588
$ws.Range("L64").Value = 8340.588
$ws.Range("N64").Value = -8836.588
$ws.Range("H67").Value = 7665.4546
$ws.Range("J67").Value = 8340.588
$ws.Range("L67").Value = 8340.588
$ws.Range("N67").Value = -10056.588
$ws.Range("H99").Value = 306.66666
$ws.Range("I99").Value = 210
$ws.Range("K99").Value = 630
$ws.Range("M99").Value = 868
$ws.Range("H106").Value = 3330.2
$ws.Range("I106").Value = 3330.2
$ws.Range("K106").Value = 3330.2
$ws.Range("M106").Value = -2699.2
$ws.Range("H129").Value = 1094.8235
$ws.Range("I129").Value = 975.75
$ws.Range("K129").Value = 2927.25
$ws.Range("M129").Value = 2072.75
$ws.Range("H135").Value = 1554.125
$ws.Range("I135").Value = 1187.1
$ws.Range("J135").Value = 2165.8333
$ws.Range("K135").Value = 10683.9
$ws.Range("L135").Value = 19492.4997
$ws.Range("M135").Value = -8148.9
$ws.Range("N135").Value = -24562.4997
$ws.Range("H138").Value = 3309.0193
$ws.Range("I138").Value = 2645.2727
$ws.Range("J138").Value = 3487.0977
$ws.Range("K138").Value = 7935.8181
$ws.Range("L138").Value = 10461.2931
$ws.Range("M138").Value = -2795.8181
$ws.Range("N138").Value = -20741.2931

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3143.013
$ws.Range("I32").Value = 2200.5374
$ws.Range("K32").Value = 2200.5374
$ws.Range("M32").Value = -1913.5374
$ws.Range("H43").Value = 31883.666
$ws.Range("J43").Value = 31883.666
$ws.Range("L43").Value = 31883.666
$ws.Range("N43").Value = -32509.666
$ws.Range("H61").Value = 1946.2285
$ws.Range("I61").Value = 1769.8125
$ws.Range("K61").Value = 1769.8125
$ws.Range("M61").Value = -1557.8125
$ws.Range("H132").Value = 2710.2964
$ws.Range("J132").Value = 3495.5386
$ws.Range("L132").Value = 10486.6158
$ws.Range("N132").Value = -15546.6158
$ws.Range("H136").Value = 1946.2285
$ws.Range("I136").Value = 1769.8125
$ws.Range("K136").Value = 5309.4375
$ws.Range("M136").Value = -2759.4375
$ws.Range("H138").Value = 59995
$ws.Range("J138").Value = 59995
$ws.Range("L138").Value = 59995
$ws.Range("N138").Value = -70275
$ws.Range("H139").Value = 85872
$ws.Range("J139").Value = 85780.664
$ws.Range("L139").Value = 85780.664
$ws.Range("N139").Value = -96060.664
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""
$ws.Range("H141").Value = 87077.8
$ws.Range("I141").Value = 105194.5
$ws.Range("J141").Value = 75000
$ws.Range("K141").Value = 105194.5
$ws.Range("L141").Value = 75000
$ws.Range("M141").Value = -100014.5
$ws.Range("N141").Value = -85360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1189.25
$ws.Range("I64").Value = 1013
$ws.Range("K64").Value = 1013
$ws.Range("M64").Value = -788
$ws.Range("H67").Value = 1189.25
$ws.Range("I67").Value = 1013
$ws.Range("K67").Value = 1013
$ws.Range("M67").Value = -233
$ws.Range("H80").Value = 328.7619
$ws.Range("I80").Value = 255.72728
$ws.Range("J80").Value = 409.1
$ws.Range("K80").Value = 255.72728
$ws.Range("L80").Value = 409.1
$ws.Range("M80").Value = 742.2727199999999
$ws.Range("N80").Value = -2405.1
$ws.Range("H83").Value = 328.7619
$ws.Range("I83").Value = 255.72728
$ws.Range("J83").Value = 409.1
$ws.Range("K83").Value = 1278.6364
$ws.Range("L83").Value = 2045.5
$ws.Range("M83").Value = 3713.3636
$ws.Range("N83").Value = -12029.5
$ws.Range("H94").Value = 6067914.5
$ws.Range("I94").Value = 11365409
$ws.Range("K94").Value = 11365409
$ws.Range("M94").Value = -11364958
$ws.Range("H132").Value = 75123.89
$ws.Range("J132").Value = 75123.89
$ws.Range("L132").Value = 75123.89
$ws.Range("N132").Value = -85243.89

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 37046.348
$ws.Range("I31").Value = 1417.4166
$ws.Range("J31").Value = 67585.42999999999
$ws.Range("K31").Value = 1417.4166
$ws.Range("L31").Value = 67585.42999999999
$ws.Range("M31").Value = -1122.4166
$ws.Range("N31").Value = -68175.42999999999
$ws.Range("H34").Value = 37046.348
$ws.Range("I34").Value = 1417.4166
$ws.Range("J34").Value = 67585.42999999999
$ws.Range("K34").Value = 1417.4166
$ws.Range("L34").Value = 67585.42999999999
$ws.Range("M34").Value = -1215.4166
$ws.Range("N34").Value = -67989.42999999999
$ws.Range("H38").Value = 2474.2
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 2474.2
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 2474.2
$ws.Range("M38").Value = ""
$ws.Range("N38").Value = -3228.2
$ws.Range("H46").Value = 2474.2
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2474.2
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 2474.2
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = -2896.2
$ws.Range("H58").Value = 1902.2858
$ws.Range("I58").Value = 1992.909
$ws.Range("J58").Value = 1843.6471
$ws.Range("K58").Value = 1992.909
$ws.Range("L58").Value = 1843.6471
$ws.Range("M58").Value = -1789.909
$ws.Range("N58").Value = -2249.6471
$ws.Range("H62").Value = 2000
$ws.Range("I62").Value = 1500
$ws.Range("K62").Value = 1500
$ws.Range("M62").Value = -876
$ws.Range("H65").Value = 2000
$ws.Range("I65").Value = 1500
$ws.Range("K65").Value = 7500
$ws.Range("M65").Value = -4380
$ws.Range("H132").Value = 51600.832
$ws.Range("I132").Value = 1719.8182
$ws.Range("K132").Value = 5159.4546
$ws.Range("M132").Value = -2629.4546
$ws.Range("H133").Value = 50326
$ws.Range("J133").Value = 50326
$ws.Range("L133").Value = 50326
$ws.Range("N133").Value = -55386
$ws.Range("H136").Value = 1902.2858
$ws.Range("I136").Value = 1992.909
$ws.Range("J136").Value = 1843.6471
$ws.Range("K136").Value = 5978.727000000001
$ws.Range("L136").Value = 5530.9413
$ws.Range("M136").Value = -3428.727000000001
$ws.Range("N136").Value = -10630.9413

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6600
$ws.Range("I56").Value = 6600
$ws.Range("K56").Value = 6600
$ws.Range("M56").Value = -6070
$ws.Range("H104").Value = 0
$ws.Range("I104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("K104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("M104").Value = ""
$ws.Range("N104").Value = ""
$ws.Range("H131").Value = 3510.2222
$ws.Range("J131").Value = 4619
$ws.Range("L131").Value = 13857
$ws.Range("N131").Value = -23937
$ws.Range("H132").Value = 2247.6072
$ws.Range("I132").Value = 1522.3334
$ws.Range("J132").Value = 2591.158
$ws.Range("K132").Value = 13701.0006
$ws.Range("L132").Value = 23320.422
$ws.Range("M132").Value = -11171.0006
$ws.Range("N132").Value = -28380.422

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 399.36365
$ws.Range("I3").Value = 267.66666
$ws.Range("J3").Value = 448.75
$ws.Range("K3").Value = 267.66666
$ws.Range("L3").Value = 448.75
$ws.Range("M3").Value = -151.66666
$ws.Range("N3").Value = -680.75
$ws.Range("H51").Value = 84814.28999999999
$ws.Range("J51").Value = 84814.28999999999
$ws.Range("L51").Value = 84814.28999999999
$ws.Range("N51").Value = -85832.28999999999
$ws.Range("H97").Value = 1832973.8
$ws.Range("I97").Value = 1985671.5
$ws.Range("K97").Value = 1985671.5
$ws.Range("M97").Value = -1985175.5
$ws.Range("H102").Value = 6475562.5
$ws.Range("I102").Value = 8549408
$ws.Range("K102").Value = 8549408
$ws.Range("M102").Value = -8547786
$ws.Range("H132").Value = 3758.6956
$ws.Range("I132").Value = 3078.3572
$ws.Range("J132").Value = 4817
$ws.Range("K132").Value = 9235.071599999999
$ws.Range("L132").Value = 14451
$ws.Range("M132").Value = -6705.071599999999
$ws.Range("N132").Value = -19511
$ws.Range("H134").Value = 34468.25
$ws.Range("J134").Value = 34468.25
$ws.Range("L134").Value = 103404.75
$ws.Range("N134").Value = -108474.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5199.6665
$ws.Range("I46").Value = 3233.1667
$ws.Range("K46").Value = 3233.1667
$ws.Range("M46").Value = -3045.1667
$ws.Range("H68").Value = 2543.1538
$ws.Range("I68").Value = 2181.125
$ws.Range("K68").Value = 2181.125
$ws.Range("M68").Value = -1432.125
$ws.Range("H71").Value = 2543.1538
$ws.Range("I71").Value = 2181.125
$ws.Range("K71").Value = 10905.625
$ws.Range("M71").Value = -7161.625
$ws.Range("H131").Value = 111633.336
$ws.Range("J131").Value = 111633.336
$ws.Range("L131").Value = 111633.336
$ws.Range("N131").Value = -121713.336
$ws.Range("H136").Value = 107681.63
$ws.Range("I136").Value = 127472.25
$ws.Range("J136").Value = 2131.6667
$ws.Range("K136").Value = 382416.75
$ws.Range("L136").Value = 6395.000100000001
$ws.Range("M136").Value = -379866.75
$ws.Range("N136").Value = -11495.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("K14").Value = 0
$ws.Range("M14").Value = ""
$ws.Range("H100").Value = 1500.6666
$ws.Range("I100").Value = 702
$ws.Range("J100").Value = 1900
$ws.Range("K100").Value = 1404
$ws.Range("L100").Value = 3800
$ws.Range("M100").Value = -863
$ws.Range("N100").Value = -4882
$ws.Range("H107").Value = 35716504
$ws.Range("I107").Value = 45457124
$ws.Range("K107").Value = 136371372
$ws.Range("M107").Value = -136369452
$ws.Range("H113").Value = 689.4074000000001
$ws.Range("J113").Value = 894.75
$ws.Range("L113").Value = 2684.25
$ws.Range("N113").Value = -7024.25
$ws.Range("H132").Value = 77764110
$ws.Range("I132").Value = 200004340
$ws.Range("J132").Value = 1363966.9
$ws.Range("K132").Value = 600013020
$ws.Range("L132").Value = 4091900.7
$ws.Range("M132").Value = -600010490
$ws.Range("N132").Value = -4096960.7
$ws.Range("H135").Value = 61904.91
$ws.Range("J135").Value = 74325.664
$ws.Range("L135").Value = 74325.664
$ws.Range("N135").Value = -84465.664
$ws.Range("H136").Value = 3153
$ws.Range("I136").Value = 1628.8667
$ws.Range("J136").Value = 5439.2
$ws.Range("K136").Value = 4886.6001
$ws.Range("L136").Value = 16317.6
$ws.Range("M136").Value = -2336.6001
$ws.Range("N136").Value = -21417.6
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = ""
